$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Added engineering floor data (rows 39-42: room number in A, seat count in B)
$ws.Range("A39").Value = 117
$ws.Range("B39").Value = 40

$ws.Range("A40").Value = 127
$ws.Range("B40").Value = 20

$ws.Range("A41").Value = 130
$ws.Range("B41").Value = 38

$ws.Range("A42").Value = 248
$ws.Range("B42").Value = 40

# Reflect the scrolled / selected viewport state after entering the new data
$ws.Range("G33:G34").Select()
